# Update building block types / version for Metabolomics MassSpec assay template

$wb = $excel.ActiveWorkbook

# 1. Bump template version on the isa_template metadata sheet
$metaSheet = $wb.Worksheets.Item("isa_template")
$metaSheet.Range("B4").Value = "1.1.9"

# 2. Rename building block headers in the annotation table on the assay sheet
$dataSheet = $wb.Worksheets.Item("3ASY03_MetabolomicsMassSpec")

# Parameter [MS sample type] -> Characteristic [MS sample type]
$dataSheet.Range("K1").Value = "Characteristic [MS sample type]"

# Parameter [chromatography column model] -> Component [chromatography column model]
$dataSheet.Range("AC1").Value = "Component [chromatography column model]"

# Parameter [chromatography guard column model] -> Component [chromatography guard column model]
$dataSheet.Range("AF1").Value = "Component [chromatography guard column model]"
